$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 218
$ws.Range("I4").Value = 218
$ws.Range("K4").Value = 218
$ws.Range("M4").Value = -104

$ws.Range("H9").Value = 120.37037
$ws.Range("I9").Value = 91.08696
$ws.Range("J9").Value = 288.75
$ws.Range("K9").Value = 91.08696
$ws.Range("L9").Value = 288.75
$ws.Range("M9").Value = 77.91304
$ws.Range("N9").Value = -626.75

$ws.Range("H32").Value = 66678012
$ws.Range("I32").Value = 142873570
$ws.Range("J32").Value = 6899.75
$ws.Range("K32").Value = 142873570
$ws.Range("L32").Value = 6899.75
$ws.Range("M32").Value = -142873244
$ws.Range("N32").Value = -7551.75

$ws.Range("H39").Value = 150.125
$ws.Range("I39").Value = 121.57143
$ws.Range("K39").Value = 364.71429
$ws.Range("M39").Value = -68.71429000000001

$ws.Range("H98").Value = 1860.3334
$ws.Range("I98").Value = 1560.1052
$ws.Range("K98").Value = 1560.1052
$ws.Range("M98").Value = -62.10519999999997

$ws.Range("H122").Value = 1860.3334
$ws.Range("I122").Value = 1560.1052
$ws.Range("K122").Value = 4680.3156
$ws.Range("M122").Value = -2230.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2720.8948
$ws.Range("I2").Value = 2806
$ws.Range("K2").Value = 2806
$ws.Range("M2").Value = -2693

$ws.Range("H45").Value = 1609
$ws.Range("I45").Value = 1609
$ws.Range("K45").Value = 1609
$ws.Range("M45").Value = -1232

$ws.Range("H97").Value = 4966.15
$ws.Range("I97").Value = 984.5333000000001
$ws.Range("K97").Value = 984.5333000000001
$ws.Range("M97").Value = -488.5333000000001

$ws.Range("H116").Value = 2720.8948
$ws.Range("I116").Value = 2806
$ws.Range("K116").Value = 2806
$ws.Range("M116").Value = -512

$ws.Range("H122").Value = 2362.08
$ws.Range("I122").Value = 2280.15
$ws.Range("K122").Value = 6840.450000000001
$ws.Range("M122").Value = -4390.450000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2720.8948
$ws.Range("I3").Value = 2806
$ws.Range("K3").Value = 2806
$ws.Range("M3").Value = -2692

$ws.Range("H11").Value = 2202.9167
$ws.Range("I11").Value = 2837.2222
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 2837.2222
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -2697.2222
$ws.Range("N11").Value = -580

$ws.Range("H22").Value = 507.5
$ws.Range("I22").Value = 486.1111
$ws.Range("K22").Value = 486.1111
$ws.Range("M22").Value = -313.1111

$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 30000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 30000
$ws.Range("N31").Value = -30504
$ws.Range("M31").ClearContents()

$ws.Range("H37").Value = 612.5
$ws.Range("I37").Value = 483.33334
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 483.33334
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = -346.33334
$ws.Range("N37").Value = -1274

$ws.Range("H86").Value = 3947369.2
$ws.Range("I86").Value = 5152023
$ws.Range("J86").Value = 32244.75
$ws.Range("K86").Value = 5152023
$ws.Range("L86").Value = 32244.75
$ws.Range("M86").Value = -5150900
$ws.Range("N86").Value = -34490.75

$ws.Range("H89").Value = 3947369.2
$ws.Range("I89").Value = 5152023
$ws.Range("J89").Value = 32244.75
$ws.Range("K89").Value = 25760115
$ws.Range("L89").Value = 161223.75
$ws.Range("M89").Value = -25754499
$ws.Range("N89").Value = -172455.75

$ws.Range("H94").Value = 2339.6428
$ws.Range("I94").Value = 2414.8096
$ws.Range("J94").Value = 2114.1428
$ws.Range("K94").Value = 2414.8096
$ws.Range("L94").Value = 2114.1428
$ws.Range("M94").Value = -1963.8096
$ws.Range("N94").Value = -3016.1428

$ws.Range("H99").Value = 1652.4584
$ws.Range("J99").Value = 3902.6
$ws.Range("L99").Value = 3902.6
$ws.Range("N99").Value = -6898.6

$ws.Range("H105").Value = 3500.7673
$ws.Range("I105").Value = 2443.889
$ws.Range("K105").Value = 2443.889
$ws.Range("M105").Value = -696.8890000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 44804
$ws.Range("I103").Value = 44804
$ws.Range("K103").Value = 44804
$ws.Range("M103").Value = -43632

$ws.Range("H105").Value = 9416
$ws.Range("I105").Value = 9416
$ws.Range("K105").Value = 9416
$ws.Range("M105").Value = -7669

$ws.Range("H132").Value = 300
$ws.Range("I132").Value = 300
$ws.Range("K132").Value = 900
$ws.Range("M132").Value = 1630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3985.8
$ws.Range("J3").Value = 12500
$ws.Range("L3").Value = 37500
$ws.Range("N3").Value = -37724

$ws.Range("H94").Value = 5262.5
$ws.Range("I94").Value = 4011.5
$ws.Range("J94").Value = 6513.5
$ws.Range("K94").Value = 12034.5
$ws.Range("L94").Value = 19540.5
$ws.Range("M94").Value = -11358.5
$ws.Range("N94").Value = -20892.5

$ws.Range("H99").Value = 883.3333
$ws.Range("I99").Value = 883.3333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2649.9999
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -403.9998999999998
$ws.Range("M99").ClearContents()

$ws.Range("H137").Value = 15925568
$ws.Range("I137").Value = 1539.75
$ws.Range("K137").Value = 4619.25
$ws.Range("M137").Value = 480.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 34820
$ws.Range("J53").Value = 30080.334
$ws.Range("L53").Value = 30080.334
$ws.Range("N53").Value = -31342.334

$ws.Range("H57").Value = 72563.336
$ws.Range("I57").Value = 17000
$ws.Range("J57").Value = 83676
$ws.Range("K57").Value = 17000
$ws.Range("L57").Value = 83676
$ws.Range("N57").Value = -85316
$ws.Range("M57").Value = -16180

$ws.Range("H80").Value = 4171.8276
$ws.Range("I80").Value = 3286.647
$ws.Range("J80").Value = 5425.8335
$ws.Range("K80").Value = 3286.647
$ws.Range("L80").Value = 5425.8335
$ws.Range("M80").Value = -2288.647
$ws.Range("N80").Value = -7421.8335

$ws.Range("H83").Value = 4171.8276
$ws.Range("I83").Value = 3286.647
$ws.Range("J83").Value = 5425.8335
$ws.Range("K83").Value = 16433.235
$ws.Range("L83").Value = 27129.1675
$ws.Range("M83").Value = -11441.235
$ws.Range("N83").Value = -37113.1675

$ws.Range("H102").Value = 43410.918
$ws.Range("I102").Value = 1360.2222
$ws.Range("K102").Value = 1360.2222
$ws.Range("M102").Value = 261.7778000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3915.3809
$ws.Range("I22").Value = 815.8
$ws.Range("J22").Value = 4884
$ws.Range("K22").Value = 815.8
$ws.Range("L22").Value = 4884
$ws.Range("M22").Value = -520.8
$ws.Range("N22").Value = -5474

$ws.Range("H27").Value = 3915.3809
$ws.Range("I27").Value = 815.8
$ws.Range("J27").Value = 4884
$ws.Range("K27").Value = 815.8
$ws.Range("L27").Value = 4884
$ws.Range("M27").Value = -708.8
$ws.Range("N27").Value = -5098

$ws.Range("H40").Value = 7059.8
$ws.Range("I40").Value = 6603.357
$ws.Range("K40").Value = 6603.357
$ws.Range("M40").Value = -6467.357

$ws.Range("H61").Value = 22999.416
$ws.Range("I61").Value = 20443.666
$ws.Range("K61").Value = 20443.666
$ws.Range("M61").Value = -20241.666

$ws.Range("H68").Value = 6586.857
$ws.Range("J68").Value = 7805.5
$ws.Range("L68").Value = 7805.5
$ws.Range("N68").Value = -9303.5

$ws.Range("H71").Value = 6586.857
$ws.Range("J71").Value = 7805.5
$ws.Range("L71").Value = 39027.5
$ws.Range("N71").Value = -46515.5

$ws.Range("H82").Value = 5612.52
$ws.Range("J82").Value = 1521.2222
$ws.Range("L82").Value = 1521.2222
$ws.Range("N82").Value = -2243.2222

$ws.Range("H85").Value = 5612.52
$ws.Range("J85").Value = 1521.2222
$ws.Range("L85").Value = 1521.2222
$ws.Range("N85").Value = -4017.2222

$ws.Range("H113").Value = 22999.416
$ws.Range("I113").Value = 20443.666
$ws.Range("K113").Value = 20443.666
$ws.Range("M113").Value = -18273.666

$ws.Range("H122").Value = 6401.4546
$ws.Range("I122").Value = 3490.6667
$ws.Range("K122").Value = 10472.0001
$ws.Range("M122").Value = -8022.000100000001

$ws.Range("H132").Value = 9640.166999999999
$ws.Range("I132").Value = 9565
$ws.Range("K132").Value = 28695
$ws.Range("M132").Value = -26165

$ws.Range("H136").Value = 5408.6
$ws.Range("J136").Value = 4797
$ws.Range("L136").Value = 14391
$ws.Range("N136").Value = -19491

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 85000
$ws.Range("J99").Value = 85000
$ws.Range("L99").Value = 85000
$ws.Range("N99").Value = -90990

$ws.Range("H100").Value = 4187.25
$ws.Range("I100").Value = 1875.25
$ws.Range("K100").Value = 3750.5
$ws.Range("M100").Value = -3209.5

$ws.Range("H122").Value = 4650.077
$ws.Range("I122").Value = 5177.364
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 15532.092
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -13082.092
$ws.Range("N122").Value = -10150
